$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 134.11376953125
$ws.Cells.Item(3, 2).Value = 134.4760131835938
$ws.Cells.Item(4, 2).Value = 138.0943145751953
$ws.Cells.Item(5, 2).Value = 133.3759002685547
$ws.Cells.Item(6, 2).Value = 134.2959899902344
$ws.Cells.Item(7, 2).Value = 134.4599609375
$ws.Cells.Item(8, 2).Value = 136.9717864990234
$ws.Cells.Item(9, 2).Value = 134.8649749755859
$ws.Cells.Item(10, 2).Value = 135.5842742919922
$ws.Cells.Item(11, 2).Value = 133.2368469238281
$ws.Cells.Item(12, 2).Value = 132.3266906738281
$ws.Cells.Item(13, 2).Value = 136.4253387451172
$ws.Cells.Item(14, 2).Value = 138.4982757568359
$ws.Cells.Item(15, 2).Value = 147.0182952880859
$ws.Cells.Item(16, 2).Value = 155.4115600585938
$ws.Cells.Item(17, 2).Value = 190.9487762451172
$ws.Cells.Item(18, 2).Value = 182.5352020263672
$ws.Cells.Item(19, 2).Value = 189.5859527587891
$ws.Cells.Item(20, 2).Value = 182.4674835205078
$ws.Cells.Item(21, 2).Value = 184.5615386962891
$ws.Cells.Item(22, 2).Value = 185.1528167724609
$ws.Cells.Item(23, 2).Value = 181.9319763183594
$ws.Cells.Item(24, 2).Value = 179.4960632324219
$ws.Cells.Item(25, 2).Value = 180.8661499023438
$ws.Cells.Item(26, 2).Value = 181.6646575927734
$ws.Cells.Item(27, 2).Value = 180.5810852050781
$ws.Cells.Item(28, 2).Value = 185.1621551513672
$ws.Cells.Item(29, 2).Value = 173.7883605957031
$ws.Cells.Item(30, 2).Value = 178.9892730712891
$ws.Cells.Item(31, 2).Value = 180.5341644287109
$ws.Cells.Item(32, 2).Value = 190.8003234863281
$ws.Cells.Item(33, 2).Value = 210.7864685058594
$ws.Cells.Item(34, 2).Value = 208.1833953857422
$ws.Cells.Item(35, 2).Value = 247.9057769775391
$ws.Cells.Item(36, 2).Value = 265.2348327636719
$ws.Cells.Item(37, 2).Value = 265.0682373046875
$ws.Cells.Item(38, 2).Value = 244.3307800292969
$ws.Cells.Item(39, 2).Value = 231.3908843994141
$ws.Cells.Item(40, 2).Value = 197.5352172851562
$ws.Cells.Item(41, 2).Value = 187.0646209716797
$ws.Cells.Item(42, 2).Value = 174.1491546630859
$ws.Cells.Item(43, 2).Value = 153.9716339111328
$ws.Cells.Item(44, 2).Value = 156.7312316894531
$ws.Cells.Item(45, 2).Value = 135.7577209472656
$ws.Cells.Item(46, 2).Value = 140.0220642089844
$ws.Cells.Item(47, 2).Value = 123.2612991333008
$ws.Cells.Item(48, 2).Value = 136.7943115234375
$ws.Cells.Item(49, 2).Value = 123.5925903320312
